# Update gh-pages to output generated at 456a3b4
# Applies the daily-refresh numeric deltas (想去人数 / "want to go" counts)
# across the four sheets, plus one ticket-status flip on 本地生活!G3.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 559
$ws.Range("F6").Value = 14
$ws.Range("F7").Value = 52
$ws.Range("F9").Value = 370
$ws.Range("F10").Value = 4474
$ws.Range("F11").Value = 4474
$ws.Range("F15").Value = 575
$ws.Range("F16").Value = 3908
$ws.Range("F17").Value = 145
$ws.Range("F18").Value = 137
$ws.Range("F20").Value = 152
$ws.Range("F21").Value = 3290
$ws.Range("F25").Value = 2761
$ws.Range("F26").Value = 101
$ws.Range("F27").Value = 111
$ws.Range("F29").Value = 121
$ws.Range("F30").Value = 158
$ws.Range("F31").Value = 154
$ws.Range("F32").Value = 64
$ws.Range("F36").Value = 122
$ws.Range("F37").Value = 5041
$ws.Range("F38").Value = 669
$ws.Range("F39").Value = 372
$ws.Range("F40").Value = 73
$ws.Range("F42").Value = 7
$ws.Range("F43").Value = 981
$ws.Range("F44").Value = 372
$ws.Range("F46").Value = 1862
$ws.Range("F47").Value = 284
$ws.Range("F49").Value = 665
$ws.Range("F50").Value = 796

# --- 演出 (Performance) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 35
$ws.Range("F21").Value = 701

# --- 本地生活 (Local Life) sheet ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("G3").Value = "不可售"

# --- 全部类型 (All Types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 559
$ws.Range("F9").Value = 14
$ws.Range("F11").Value = 52
$ws.Range("F13").Value = 370
$ws.Range("F14").Value = 4474
$ws.Range("F15").Value = 4474
$ws.Range("F16").Value = 35
$ws.Range("F21").Value = 575
$ws.Range("F22").Value = 3908
$ws.Range("F23").Value = 145
$ws.Range("F24").Value = 137
$ws.Range("F25").Value = 3290
$ws.Range("F26").Value = 2761
$ws.Range("F27").Value = 101
$ws.Range("F28").Value = 111
$ws.Range("F29").Value = 121
$ws.Range("F30").Value = 158
$ws.Range("F31").Value = 154
$ws.Range("F32").Value = 64
$ws.Range("F36").Value = 122
$ws.Range("F39").Value = 5041
$ws.Range("F41").Value = 669
$ws.Range("F42").Value = 372
$ws.Range("F44").Value = 73
$ws.Range("F45").Value = 981
$ws.Range("F46").Value = 372
$ws.Range("F48").Value = 1862
$ws.Range("F49").Value = 284
$ws.Range("F50").Value = 665
$ws.Range("F51").Value = 796
